# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the file
# "3db582a9-0042-49a5-a2ee-b7c3e39538de" row on both the zh-cn and de-de
# status sheets, reflecting a new handback cycle that just completed.

$wb = $excel.ActiveWorkbook

# zh-cn sheet (row 2 == 3db582a9-0042-49a5-a2ee-b7c3e39538de)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 03:34:09"
$wsZhCn.Range("H2").Value = "2016-03-14 03:34:27"

# de-de sheet (row 2 == 3db582a9-0042-49a5-a2ee-b7c3e39538de)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 03:34:12"
$wsDeDe.Range("H2").Value = "2016-03-14 03:34:31"
